$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 12 ("2021年") appended after the existing last data row (11).
$ws.Range("A12").Value = "2021年"
$ws.Range("B12").Value = 14640311.8
$ws.Range("C12").Value = "'"
$ws.Range("D12").Value = 11166787.5
$ws.Range("E12").Value = 43695641.1
$ws.Range("F12").Value = 14230507.3
$ws.Range("G12").Value = 46290515.6
$ws.Range("H12").Value = 110451212.2
$ws.Range("I12").Value = 32259844.7
$ws.Range("J12").Value = 7925189.8
$ws.Range("K12").Value = 34838702.4
$ws.Range("L12").Value = "'"
$ws.Range("M12").Value = 548351401.1
$ws.Range("N12").Value = 107796846.6
$ws.Range("O12").Value = 88959483
$ws.Range("P12").Value = 10365703.3
$ws.Range("Q12").Value = 57175157.3
$ws.Range("R12").Value = "'"
$ws.Range("S12").Value = "'"
$ws.Range("T12").Value = 1004563.1
$ws.Range("U12").Value = 817425737.5

# Copy the formatting from row 11 down onto the new row so the year label in
# column A keeps the bold/bordered/centered header-style used by every other
# "year" cell in column A, and the empty text cells (C/L/R/S, used throughout
# this sheet for categories with no data for a given year) keep the plain,
# unstyled look of their counterparts in the row above.
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("C11").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("L11").Copy()
$ws.Range("L12").PasteSpecial(-4122)
$ws.Range("R11").Copy()
$ws.Range("R12").PasteSpecial(-4122)
$ws.Range("S11").Copy()
$ws.Range("S12").PasteSpecial(-4122)
